# Update countries & provincias Spain
# Refreshes the COVID-19 "Pais" sheet with newer case counts and re-applies
# the resulting country ranking (some countries swap rows because their
# "Casos totales" changed relative to their neighbour), plus bumps the
# "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- "Datos actualizados" timestamp (row 1) ----
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 17 de Octubre de 2020 a las 17:29"

# ---- Row 4: Estados Unidos (new totals) ----
$ws.Cells.Item(4, 2).Value = 8304539
$ws.Cells.Item(4, 3).Value = 16261
$ws.Cells.Item(4, 4).Value = 5402456
$ws.Cells.Item(4, 5).Value = 2678207
$ws.Cells.Item(4, 7).Value = 232
$ws.Cells.Item(4, 8).Value = 223876

# ---- Row 5: India (new totals) ----
$ws.Cells.Item(5, 2).Value = 7443233
$ws.Cells.Item(5, 3).Value = 12598
$ws.Cells.Item(5, 4).Value = 6534590
$ws.Cells.Item(5, 5).Value = 795471
$ws.Cells.Item(5, 7).Value = 140
$ws.Cells.Item(5, 8).Value = 113172

# ---- Rows 14-15: Reino Unido overtakes Sudafrica ----
$ws.Cells.Item(14, 1).Value = "Reino Unido"
$ws.Cells.Item(14, 2).Value = 705428
$ws.Cells.Item(14, 3).Value = 16171
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 7).Value = 150
$ws.Cells.Item(14, 8).Value = 43579

$ws.Cells.Item(15, 1).Value = "Sudafrica"
$ws.Cells.Item(15, 2).Value = 700203
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 629260
$ws.Cells.Item(15, 5).Value = 52573
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 18370

# ---- Row 17: Chile (new totals) ----
$ws.Cells.Item(17, 2).Value = 490003
$ws.Cells.Item(17, 3).Value = 1813
$ws.Cells.Item(17, 4).Value = 462712
$ws.Cells.Item(17, 5).Value = 13703
$ws.Cells.Item(17, 7).Value = 59
$ws.Cells.Item(17, 8).Value = 13588

# ---- Row 19: Italia (new totals) ----
$ws.Cells.Item(19, 2).Value = 402536
$ws.Cells.Item(19, 3).Value = 10925
$ws.Cells.Item(19, 4).Value = 249127
$ws.Cells.Item(19, 5).Value = 116935
$ws.Cells.Item(19, 7).Value = 47
$ws.Cells.Item(19, 8).Value = 36474

# ---- Rows 21-22: Alemania overtakes Indonesia ----
$ws.Cells.Item(21, 1).Value = "Alemania"
$ws.Cells.Item(21, 2).Value = 358510
$ws.Cells.Item(21, 3).Value = 1718
$ws.Cells.Item(21, 4).Value = 290000
$ws.Cells.Item(21, 5).Value = 58672
$ws.Cells.Item(21, 7).Value = 2
$ws.Cells.Item(21, 8).Value = 9838

$ws.Cells.Item(22, 1).Value = "Indonesia"
$ws.Cells.Item(22, 2).Value = 357762
$ws.Cells.Item(22, 3).Value = 4301
$ws.Cells.Item(22, 4).Value = 281592
$ws.Cells.Item(22, 5).Value = 63739
$ws.Cells.Item(22, 7).Value = 84
$ws.Cells.Item(22, 8).Value = 12431

# ---- Row 41: Republica Dominicana (new totals) ----
$ws.Cells.Item(41, 2).Value = 120925
$ws.Cells.Item(41, 3).Value = 475
$ws.Cells.Item(41, 4).Value = 97575
$ws.Cells.Item(41, 5).Value = 21155
$ws.Cells.Item(41, 7).Value = 3
$ws.Cells.Item(41, 8).Value = 2195

# ---- Row 48: Guatemala (new totals) ----
$ws.Cells.Item(48, 2).Value = 101028
$ws.Cells.Item(48, 3).Value = 597
$ws.Cells.Item(48, 4).Value = 90001
$ws.Cells.Item(48, 5).Value = 7512
$ws.Cells.Item(48, 7).Value = 37
$ws.Cells.Item(48, 8).Value = 3515

# ---- Row 51: Japon (new totals) ----
$ws.Cells.Item(51, 2).Value = 92063
$ws.Cells.Item(51, 3).Value = 632
$ws.Cells.Item(51, 4).Value = 85030
$ws.Cells.Item(51, 5).Value = 5372
$ws.Cells.Item(51, 7).Value = 11
$ws.Cells.Item(51, 8).Value = 1661

# ---- Row 59: Moldavia (new totals) ----
$ws.Cells.Item(59, 2).Value = 66652
$ws.Cells.Item(59, 3).Value = 792
$ws.Cells.Item(59, 4).Value = 47230
$ws.Cells.Item(59, 5).Value = 17853
$ws.Cells.Item(59, 7).Value = 20
$ws.Cells.Item(59, 8).Value = 1569

# ---- Row 65: Singapur (new totals) ----
$ws.Cells.Item(65, 4).Value = 57798
$ws.Cells.Item(65, 5).Value = 78

# ---- Row 95: Albania (new totals) ----
$ws.Cells.Item(95, 2).Value = 16774
$ws.Cells.Item(95, 3).Value = 273
$ws.Cells.Item(95, 4).Value = 10001
$ws.Cells.Item(95, 5).Value = 6325
$ws.Cells.Item(95, 7).Value = 5
$ws.Cells.Item(95, 8).Value = 448

# ---- Row 115: Jamaica (new totals) ----
$ws.Cells.Item(115, 2).Value = 8195
$ws.Cells.Item(115, 3).Value = 63
$ws.Cells.Item(115, 5).Value = 4374
$ws.Cells.Item(115, 7).Value = 6
$ws.Cells.Item(115, 8).Value = 168

# ---- Row 122: Cuba (new totals) ----
$ws.Cells.Item(122, 2).Value = 6170
$ws.Cells.Item(122, 3).Value = 52
$ws.Cells.Item(122, 4).Value = 5753
$ws.Cells.Item(122, 7).Value = 1
$ws.Cells.Item(122, 8).Value = 125

# ---- Row 123: Malaui (new totals) ----
$ws.Cells.Item(123, 2).Value = 5852
$ws.Cells.Item(123, 3).Value = 10
$ws.Cells.Item(123, 4).Value = 4740
$ws.Cells.Item(123, 5).Value = 931

# ---- Row 127: Republica de Yibuti (new totals) ----
$ws.Cells.Item(127, 2).Value = 5452
$ws.Cells.Item(127, 3).Value = 3
$ws.Cells.Item(127, 4).Value = 5373
$ws.Cells.Item(127, 5).Value = 18

# ---- Row 183: Eritrea (new totals) ----
$ws.Cells.Item(183, 2).Value = 452
$ws.Cells.Item(183, 3).Value = 30
$ws.Cells.Item(183, 4).Value = 388
$ws.Cells.Item(183, 5).Value = 64

# ---- Row 192: Liechtenstein (new totals) ----
$ws.Cells.Item(192, 2).Value = 217
$ws.Cells.Item(192, 3).Value = 25
$ws.Cells.Item(192, 5).Value = 84

# ---- Rows 216-217: Montserrat / Islas Malvinas swap places ----
$ws.Cells.Item(216, 1).Value = "Islas Malvinas"
$ws.Cells.Item(216, 4).Value = 13
$ws.Cells.Item(216, 8).Value = 0

$ws.Cells.Item(217, 1).Value = "Montserrat"
$ws.Cells.Item(217, 4).Value = 12
$ws.Cells.Item(217, 8).Value = 1
